$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 1.993530333333333
$ws.Cells.Item(2, 8).Value = 5.980591
$ws.Cells.Item(2, 9).Value = 0.2657535748877558
$ws.Cells.Item(2, 10).Value = 0.2657535748877558
$ws.Cells.Item(2, 13).Value = 1.123319
$ws.Cells.Item(2, 14).Value = 3.369957
$ws.Cells.Item(2, 15).Value = 0.05053686506648315
$ws.Cells.Item(2, 16).Value = 0.05053686506648315
$ws.Cells.Item(2, 17).Value = 2.239370500509667
$ws.Cells.Item(2, 18).Value = 20.154334504587
$ws.Cells.Item(2, 19).Value = 0.01343035255503804
$ws.Cells.Item(2, 20).Value = 0.01343035255503804

$ws.Cells.Item(3, 7).Value = 1.993530333333333
$ws.Cells.Item(3, 8).Value = 5.980591
$ws.Cells.Item(3, 9).Value = 0.2657535748877558
$ws.Cells.Item(3, 10).Value = 0.2657535748877558
$ws.Cells.Item(3, 15).Value = 0.5042195746532222
$ws.Cells.Item(3, 16).Value = 0.5042195746532223
$ws.Cells.Item(3, 17).Value = 22.34278758234289
$ws.Cells.Item(3, 18).Value = 201.085088241086
$ws.Cells.Item(3, 19).Value = 0.1339981544924775
$ws.Cells.Item(3, 20).Value = 0.1339981544924775

$ws.Cells.Item(4, 7).Value = 1.993530333333333
$ws.Cells.Item(4, 8).Value = 5.980591
$ws.Cells.Item(4, 9).Value = 0.2657535748877558
$ws.Cells.Item(4, 10).Value = 0.2657535748877558
$ws.Cells.Item(4, 13).Value = 4.958620666666667
$ws.Cells.Item(4, 14).Value = 14.875862
$ws.Cells.Item(4, 15).Value = 0.2230827962023326
$ws.Cells.Item(4, 16).Value = 0.2230827962023326
$ws.Cells.Item(4, 17).Value = 9.885160710493556
$ws.Cells.Item(4, 18).Value = 88.966446394442
$ws.Cells.Item(4, 19).Value = 0.05928505058672656
$ws.Cells.Item(4, 20).Value = 0.05928505058672655

$ws.Cells.Item(5, 7).Value = 1.993530333333333
$ws.Cells.Item(5, 8).Value = 5.980591
$ws.Cells.Item(5, 9).Value = 0.2657535748877558
$ws.Cells.Item(5, 10).Value = 0.2657535748877558
$ws.Cells.Item(5, 13).Value = 4.938126
$ws.Cells.Item(5, 14).Value = 14.814378
$ws.Cells.Item(5, 15).Value = 0.222160764077962
$ws.Cells.Item(5, 16).Value = 0.222160764077962
$ws.Cells.Item(5, 17).Value = 9.844303970822001
$ws.Cells.Item(5, 18).Value = 88.598735737398
$ws.Cells.Item(5, 19).Value = 0.05904001725351371
$ws.Cells.Item(5, 20).Value = 0.05904001725351371

$ws.Cells.Item(6, 9).Value = 0.324410035374171
$ws.Cells.Item(6, 10).Value = 0.324410035374171
$ws.Cells.Item(6, 13).Value = 1.123319
$ws.Cells.Item(6, 14).Value = 3.369957
$ws.Cells.Item(6, 15).Value = 0.05053686506648315
$ws.Cells.Item(6, 16).Value = 0.05053686506648315
$ws.Cells.Item(6, 17).Value = 2.733638723742667
$ws.Cells.Item(6, 18).Value = 24.602748513684
$ws.Cells.Item(6, 19).Value = 0.0163946661839175
$ws.Cells.Item(6, 20).Value = 0.0163946661839175

$ws.Cells.Item(7, 9).Value = 0.324410035374171
$ws.Cells.Item(7, 10).Value = 0.324410035374171
$ws.Cells.Item(7, 15).Value = 0.5042195746532222
$ws.Cells.Item(7, 16).Value = 0.5042195746532223
$ws.Cells.Item(7, 19).Value = 0.1635738900496013
$ws.Cells.Item(7, 20).Value = 0.1635738900496013

$ws.Cells.Item(8, 9).Value = 0.324410035374171
$ws.Cells.Item(8, 10).Value = 0.324410035374171
$ws.Cells.Item(8, 13).Value = 4.958620666666667
$ws.Cells.Item(8, 14).Value = 14.875862
$ws.Cells.Item(8, 15).Value = 0.2230827962023326
$ws.Cells.Item(8, 16).Value = 0.2230827962023326
$ws.Cells.Item(8, 17).Value = 12.06698851417156
$ws.Cells.Item(8, 18).Value = 108.602896627544
$ws.Cells.Item(8, 19).Value = 0.0723702978073677
$ws.Cells.Item(8, 20).Value = 0.07237029780736769

$ws.Cells.Item(9, 9).Value = 0.324410035374171
$ws.Cells.Item(9, 10).Value = 0.324410035374171
$ws.Cells.Item(9, 13).Value = 4.938126
$ws.Cells.Item(9, 14).Value = 14.814378
$ws.Cells.Item(9, 15).Value = 0.222160764077962
$ws.Cells.Item(9, 16).Value = 0.222160764077962
$ws.Cells.Item(9, 17).Value = 12.017113977704
$ws.Cells.Item(9, 18).Value = 108.154025799336
$ws.Cells.Item(9, 19).Value = 0.0720711813332845
$ws.Cells.Item(9, 20).Value = 0.0720711813332845

$ws.Cells.Item(10, 7).Value = 2.135898
$ws.Cells.Item(10, 8).Value = 6.407693999999999
$ws.Cells.Item(10, 9).Value = 0.2847323261675683
$ws.Cells.Item(10, 10).Value = 0.2847323261675683
$ws.Cells.Item(10, 13).Value = 1.123319
$ws.Cells.Item(10, 14).Value = 3.369957
$ws.Cells.Item(10, 15).Value = 0.05053686506648315
$ws.Cells.Item(10, 16).Value = 0.05053686506648315
$ws.Cells.Item(10, 17).Value = 2.399294805462
$ws.Cells.Item(10, 18).Value = 21.593653249158
$ws.Cells.Item(10, 19).Value = 0.01438947914759627
$ws.Cells.Item(10, 20).Value = 0.01438947914759627

$ws.Cells.Item(11, 7).Value = 2.135898
$ws.Cells.Item(11, 8).Value = 6.407693999999999
$ws.Cells.Item(11, 9).Value = 0.2847323261675683
$ws.Cells.Item(11, 10).Value = 0.2847323261675683
$ws.Cells.Item(11, 15).Value = 0.5042195746532222
$ws.Cells.Item(11, 16).Value = 0.5042195746532223
$ws.Cells.Item(11, 17).Value = 23.93839437183599
$ws.Cells.Item(11, 18).Value = 215.445549346524
$ws.Cells.Item(11, 19).Value = 0.1435676123902338
$ws.Cells.Item(11, 20).Value = 0.1435676123902338

$ws.Cells.Item(12, 7).Value = 2.135898
$ws.Cells.Item(12, 8).Value = 6.407693999999999
$ws.Cells.Item(12, 9).Value = 0.2847323261675683
$ws.Cells.Item(12, 10).Value = 0.2847323261675683
$ws.Cells.Item(12, 13).Value = 4.958620666666667
$ws.Cells.Item(12, 14).Value = 14.875862
$ws.Cells.Item(12, 15).Value = 0.2230827962023326
$ws.Cells.Item(12, 16).Value = 0.2230827962023326
$ws.Cells.Item(12, 17).Value = 10.591107964692
$ws.Cells.Item(12, 18).Value = 95.319971682228
$ws.Cells.Item(12, 19).Value = 0.06351888349065572
$ws.Cells.Item(12, 20).Value = 0.06351888349065572

$ws.Cells.Item(13, 7).Value = 2.135898
$ws.Cells.Item(13, 8).Value = 6.407693999999999
$ws.Cells.Item(13, 9).Value = 0.2847323261675683
$ws.Cells.Item(13, 10).Value = 0.2847323261675683
$ws.Cells.Item(13, 13).Value = 4.938126
$ws.Cells.Item(13, 14).Value = 14.814378
$ws.Cells.Item(13, 15).Value = 0.222160764077962
$ws.Cells.Item(13, 16).Value = 0.222160764077962
$ws.Cells.Item(13, 17).Value = 10.547333447148
$ws.Cells.Item(13, 18).Value = 94.926001024332
$ws.Cells.Item(13, 19).Value = 0.06325635113908244
$ws.Cells.Item(13, 20).Value = 0.06325635113908246

$ws.Cells.Item(14, 7).Value = 0.9384586666666667
$ws.Cells.Item(14, 8).Value = 2.815376
$ws.Cells.Item(14, 9).Value = 0.125104063570505
$ws.Cells.Item(14, 10).Value = 0.125104063570505
$ws.Cells.Item(14, 13).Value = 1.123319
$ws.Cells.Item(14, 14).Value = 3.369957
$ws.Cells.Item(14, 15).Value = 0.05053686506648315
$ws.Cells.Item(14, 16).Value = 0.05053686506648315
$ws.Cells.Item(14, 17).Value = 1.054188450981334
$ws.Cells.Item(14, 18).Value = 9.487696058832002
$ws.Cells.Item(14, 19).Value = 0.006322367179931343
$ws.Cells.Item(14, 20).Value = 0.006322367179931343

$ws.Cells.Item(15, 7).Value = 0.9384586666666667
$ws.Cells.Item(15, 8).Value = 2.815376
$ws.Cells.Item(15, 9).Value = 0.125104063570505
$ws.Cells.Item(15, 10).Value = 0.125104063570505
$ws.Cells.Item(15, 15).Value = 0.5042195746532222
$ws.Cells.Item(15, 16).Value = 0.5042195746532223
$ws.Cells.Item(15, 17).Value = 10.51791502418844
$ws.Cells.Item(15, 18).Value = 94.661235217696
$ws.Cells.Item(15, 19).Value = 0.06307991772090972
$ws.Cells.Item(15, 20).Value = 0.06307991772090973

$ws.Cells.Item(16, 7).Value = 0.9384586666666667
$ws.Cells.Item(16, 8).Value = 2.815376
$ws.Cells.Item(16, 9).Value = 0.125104063570505
$ws.Cells.Item(16, 10).Value = 0.125104063570505
$ws.Cells.Item(16, 13).Value = 4.958620666666667
$ws.Cells.Item(16, 14).Value = 14.875862
$ws.Cells.Item(16, 15).Value = 0.2230827962023326
$ws.Cells.Item(16, 16).Value = 0.2230827962023326
$ws.Cells.Item(16, 17).Value = 4.653460539345779
$ws.Cells.Item(16, 18).Value = 41.88114485411201
$ws.Cells.Item(16, 19).Value = 0.02790856431758264
$ws.Cells.Item(16, 20).Value = 0.02790856431758264

$ws.Cells.Item(17, 7).Value = 0.9384586666666667
$ws.Cells.Item(17, 8).Value = 2.815376
$ws.Cells.Item(17, 9).Value = 0.125104063570505
$ws.Cells.Item(17, 10).Value = 0.125104063570505
$ws.Cells.Item(17, 13).Value = 4.938126
$ws.Cells.Item(17, 14).Value = 14.814378
$ws.Cells.Item(17, 15).Value = 0.222160764077962
$ws.Cells.Item(17, 16).Value = 0.222160764077962
$ws.Cells.Item(17, 17).Value = 4.634227141792
$ws.Cells.Item(17, 18).Value = 41.70804427612801
$ws.Cells.Item(17, 19).Value = 0.02779321435208133
$ws.Cells.Item(17, 20).Value = 0.02779321435208133
